$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edit: TIMEPRD columns J (09:00) and K (10:00) go from 0 to 1 for every day row ---
foreach ($r in 4..10) {
    $ws.Cells.Item($r, 10).Value = 1   # column J
    $ws.Cells.Item($r, 11).Value = 1   # column K
}

# --- Re-apply the header formatting to the merged title row (B1:Y1) ---
$headerRange = $ws.Range("B1:Y1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop

# --- Selection left where the author last clicked ---
$ws.Range("O9").Select() | Out-Null

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1              # xlPortrait

# --- Footer: corporate sensitivity-label footer inserted on save ---
$ws.PageSetup.LeftFooter = "&1#&`"Calibri`"&10&K0000FFGizlilik Sınıflandırması : Hizmete Özel"
